# Insert a new match row (AZERBAIJAN - PREMIER LEAGUE: Zira vs Araz) above the
# current row 2 (INDIA - ISL: Goa vs Punjab), pushing the existing rows 2 and 3
# down to rows 3 and 4. Also refresh the odds on the INDIA - ISL row, which
# moved to row 3, since they were updated alongside the insertion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 2:3 down to 3:4, making room for the new row 2.
$ws.Rows(2).Insert()

# New row 2 values: AZERBAIJAN - PREMIER LEAGUE, Zira vs Araz.
# The date cell is prefixed with a literal-text apostrophe so Excel keeps it
# as the string "06/11/2024" instead of auto-converting it to a date serial
# (matching every other Date cell in the column, which are plain text).
$newRow = "OxKyCg3k","'06/11/2024","12:00","AZERBAIJAN - PREMIER LEAGUE","Zira","Araz",1.98,2.77,4.5,2.67,1.82,5.2,1.14,4.3,1.6,2.07,2.67,1.36,1.6,2.07,2.25,1.5,4.85,7.8,9.25,18,21,50,4.6,5.8,22,175,700,8.25,23,16.5,90,70,90,3.55,10.75,25,45,110,450,2.02,8.25,110,5.9,30,45,200,300,450

for ($i = 0; $i -lt $newRow.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $newRow[$i]
}

# The inserted row picked up bold/boxed formatting copied from the row above
# (the header row), and the leading apostrophe above leaves a "number stored
# as text" quote-prefix marker. Clear all of that now that the values are in
# place, so the row ends up with the plain (unstyled) look used by every
# other data row.
$ws.Rows(2).ClearFormats()

# The INDIA - ISL match (now on row 3 after the shift) has several odds that
# were updated at the same time as the new match was added.
$row3Updates = @{
    7  = 1.87   # G3  Odd_H_FT
    10 = 2.6    # J3  Odd_H_HT
    11 = 2.25   # K3  Odd_D_HT
    15 = 1.25   # O3  Odd_Over15_FT
    16 = 3.75   # P3  Odd_Under15_FT
    17 = 1.83   # Q3  Odd_Over25_FT
    18 = 1.98   # R3  Odd_Under25_FT
    19 = 1.36   # S3  Odd_Over05_HT
    20 = 3      # T3  Odd_Under05_HT
    21 = 1.73   # U3  Odd_BTTS_Yes
    22 = 2      # V3  Odd_BTTS_No
    23 = 8      # W3  Odd_CS_1-0
    29 = 12     # AC3 Odd_CS_0-0
    33 = 201    # AG3 Odd_CS_4-4
    34 = 12     # AH3 Odd_CS_0-1
    35 = 19     # AI3 Odd_CS_0-2
    38 = 29     # AL3 Odd_CS_1-3
    39 = 34     # AM3 Odd_CS_2-3
    45 = 126    # AS3 Odd_CS_3-2_HT
    46 = 3      # AT3 Odd_CS_0-0_HT
    48 = 51     # AV3 Odd_CS_2-2_HT
    51 = 26     # AY3 Odd_CS_1-2_HT
    52 = 67     # AZ3 Odd_CS_0-3_HT
    53 = 81     # BA3 Odd_CS_1-3_HT
}

foreach ($col in $row3Updates.Keys) {
    $ws.Cells.Item(3, $col).Value = $row3Updates[$col]
}
